# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp.
# - Peru overtakes Colombia in total cases -> rows 8/9 swap countries,
#   row 8 gets fresh Peru figures, row 9 keeps Colombia's previous figures.
# - Birmania overtakes Malaui -> rows 113/114 swap the same way.
# - Santa Lucia overtakes Timor Oriental -> rows 204/205 swap countries.
# - Montserrat overtakes Islas Malvinas -> rows 214/215 swap countries.
# - A handful of other countries simply got refreshed case counts
#   (Lesoto row 162, Butan row 188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 06:14"

# --- Peru / Colombia (rows 8-9) ---
$ws.Range("A8").Value = "Peru"
$ws.Range("B8").Value = 768895
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 615255
$ws.Range("E8").Value = 122271
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 31369

$ws.Range("A9").Value = "Colombia"
$ws.Range("B9").Value = 765076
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 633199
$ws.Range("E9").Value = 107669
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 24208

# --- Birmania / Malaui (rows 113-114) ---
$ws.Range("A113").Value = "Birmania"
$ws.Range("B113").Value = 5805
$ws.Range("C113").Value = 264
$ws.Range("D113").Value = 1260
$ws.Range("E113").Value = 4451
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 94

$ws.Range("A114").Value = "Malaui"
$ws.Range("B114").Value = 5731
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 4040
$ws.Range("E114").Value = 1512
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 179

# --- Lesoto refresh (row 162) ---
$ws.Range("B162").Value = 1424
$ws.Range("E162").Value = 637

# --- Butan refresh (row 188) ---
$ws.Range("B188").Value = 261
$ws.Range("C188").Value = 2
$ws.Range("D188").Value = 192

# --- Santa Lucia / Timor Oriental (rows 204-205) ---
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Timor Oriental"

# --- Montserrat / Islas Malvinas (rows 214-215) ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
